# Add an "Outliers" note after the existing "Missing values" note.
#
# Before:
#   <one paragraph> "Missing values: df.null "  (bookmark "_GoBack" wraps "Missing values")
#
# After:
#   <paragraph 1> "Missing values: df.null "                         (no bookmark)
#   <paragraph 2> "Outliers: Are the data points that are different
#                  from the other points, it will cause measurement
#                  error, data entry error -> use boxplot, or IQR"   (bookmark "_GoBack", trailing space)

$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "Missing values" note.
# Walk every paragraph and remember the last match, mirroring the fact
# that this note sits at the very end of the document.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Missing values*") {
        $target = $p
    }
}

if ($target -eq $null) {
    $target = $d.Paragraphs.Last
}

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
         '<w:r><w:t xml:space="preserve">Missing values: </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>df.null</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '</w:p>' +
       '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
         '<w:r><w:t xml:space="preserve">Outliers: </w:t></w:r>' +
         '<w:r><w:t>Are the data points that are different from the other points</w:t></w:r>' +
         '<w:r><w:t>, it will cause measurement error, data entry error</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> -&gt; use boxplot, or IQR</w:t></w:r>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
         '<w:bookmarkEnd w:id="0"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '</w:p>'

# Replacing the whole paragraph range (text + trailing paragraph mark) with
# the two-paragraph XML fragment splits it into the two final paragraphs in
# one shot, without leaving a stray empty paragraph behind.
[void]$target.Range.InsertXML($xml)
